$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New player stat blocks (rows 29-40) appended after the existing data (which ended at row 28).
# Each player occupies 3 rows: Group1, Group2, Difference.
# Fill style alternates green/yellow per player block, continuing the existing pattern
# (row 28 was green, so the next block - Brian O'Neill - is yellow).

# --- Brian O'Neill (rows 29-31) : yellow block ---
$ws.Range("A29").Value = "Brian O'Neill"
$ws.Range("B29").Value = "Group1"
$ws.Range("C29").Value = 0.3333333333333333
$ws.Range("D29").Value = 0.3333333333333333
$ws.Range("E29").Value = 0

$ws.Range("A30").Value = "Brian O'Neill"
$ws.Range("B30").Value = "Group2"
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0

$ws.Range("A31").Value = "Brian O'Neill"
$ws.Range("B31").Value = "Difference"
$ws.Range("C31").Value = -0.3333333333333333
$ws.Range("D31").Value = -0.3333333333333333
$ws.Range("E31").Value = 0

# --- Garrett Bradbury (rows 32-34) : green block, no stats available ---
$ws.Range("A32").Value = "Garrett Bradbury"
$ws.Range("B32").Value = "Group1"

$ws.Range("A33").Value = "Garrett Bradbury"
$ws.Range("B33").Value = "Group2"

$ws.Range("A34").Value = "Garrett Bradbury"
$ws.Range("B34").Value = "Difference"

# --- Erik McCoy (rows 35-37) : yellow block ---
$ws.Range("A35").Value = "Erik McCoy"
$ws.Range("B35").Value = "Group1"
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0

$ws.Range("A36").Value = "Erik McCoy"
$ws.Range("B36").Value = "Group2"
$ws.Range("C36").Value = 0.6666666666666666
$ws.Range("D36").Value = 0.3333333333333333
$ws.Range("E36").Value = 0.3333333333333333

$ws.Range("A37").Value = "Erik McCoy"
$ws.Range("B37").Value = "Difference"
$ws.Range("C37").Value = 0.6666666666666666
$ws.Range("D37").Value = 0.3333333333333333
$ws.Range("E37").Value = 0.3333333333333333

# --- Isaiah Wynn (rows 38-40) : green block ---
$ws.Range("A38").Value = "Isaiah Wynn"
$ws.Range("B38").Value = "Group1"
$ws.Range("C38").Value = 0.3333333333333333
$ws.Range("D38").Value = 0.3333333333333333
$ws.Range("E38").Value = 0

$ws.Range("A39").Value = "Isaiah Wynn"
$ws.Range("B39").Value = "Group2"
$ws.Range("C39").Value = 0.3333333333333333
$ws.Range("D39").Value = 0.3333333333333333
$ws.Range("E39").Value = 0

$ws.Range("A40").Value = "Isaiah Wynn"
$ws.Range("B40").Value = "Difference"
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 0

# Apply the alternating fill formatting to match the rest of the sheet.
# Existing yellow (style s=3) blocks are e.g. rows 23-25 (Tytus Howard); existing
# green (style s=2) blocks are e.g. rows 26-28 (Wyatt Teller).
$ws.Range("A23:E25").Copy()
$ws.Range("A29:E31").PasteSpecial(-4122)

$ws.Range("A26:E28").Copy()
$ws.Range("A32:E34").PasteSpecial(-4122)

$ws.Range("A23:E25").Copy()
$ws.Range("A35:E37").PasteSpecial(-4122)

$ws.Range("A26:E28").Copy()
$ws.Range("A38:E40").PasteSpecial(-4122)

$excel.CutCopyMode = 0
